# Updated cryptos list on Fri May 31 03:46:40 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the value to be written as text (matching the original inlineStr
    # cells) instead of letting Excel auto-convert numeric-looking strings,
    # then restore the default "Normal" style so no stray number-format is
    # left behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "68.574.49"
Set-TextValue $ws.Range("E2") "  +0.78%  "

Set-TextValue $ws.Range("D3") "3.763.48"
Set-TextValue $ws.Range("E3") "  -0.59%  "

Set-TextValue $ws.Range("E4") "  -0.03%  "

Set-TextValue $ws.Range("D5") "594.18"
Set-TextValue $ws.Range("E5") "  -0.53%  "

Set-TextValue $ws.Range("D6") "167.37"
Set-TextValue $ws.Range("E6") "  -1.48%  "

Set-TextValue $ws.Range("D7") "3.762.29"
Set-TextValue $ws.Range("E7") "  -0.60%  "

Set-TextValue $ws.Range("E8") "  -0.07%  "

Set-TextValue $ws.Range("D9") "0.522"
Set-TextValue $ws.Range("E9") "  -0.94%  "

Set-TextValue $ws.Range("E10") "  -2.80%  "

Set-TextValue $ws.Range("D11") "6.42"
Set-TextValue $ws.Range("E11") "  -1.44%  "

Set-TextValue $ws.Range("E12") "  -1.02%  "

Set-TextValue $ws.Range("D13") "0.0000260"
Set-TextValue $ws.Range("E13") "  -6.88%  "

Set-TextValue $ws.Range("D14") "36.18"
Set-TextValue $ws.Range("E14") "  -1.43%  "

Set-TextValue $ws.Range("D15") "4.396.24"
Set-TextValue $ws.Range("E15") "  -0.52%  "

Set-TextValue $ws.Range("D16") "3.764.50"
Set-TextValue $ws.Range("E16") "  -0.56%  "

Set-TextValue $ws.Range("D17") "68.535.60"
Set-TextValue $ws.Range("E17") "  +0.86%  "

Set-TextValue $ws.Range("D18") "17.93"
Set-TextValue $ws.Range("E18") "  -4.48%  "

Set-TextValue $ws.Range("E19") "  +0.77%  "

Set-TextValue $ws.Range("E20") "  -2.87%  "

Set-TextValue $ws.Range("D21") "10.74"
Set-TextValue $ws.Range("E21") "  +1.34%  "

Set-TextValue $ws.Range("D22") "465.48"
Set-TextValue $ws.Range("E22") "  -0.60%  "

Set-TextValue $ws.Range("E23") "  -3.23%  "

# Rows 24 and 25 swapped places (Litecoin <-> PEPE) with refreshed prices
Set-TextValue $ws.Range("B24") "PEPE"
Set-TextValue $ws.Range("C24") "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws.Range("D24") "0.0000148"
Set-TextValue $ws.Range("E24") "  -1.65%  "

Set-TextValue $ws.Range("B25") "Litecoin"
Set-TextValue $ws.Range("C25") "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D25") "84.18"
Set-TextValue $ws.Range("E25") "  +0.49%  "

Set-TextValue $ws.Range("E26") "  -2.81%  "

Set-TextValue $ws.Range("D27") "11.98"
Set-TextValue $ws.Range("E27") "  -1.56%  "

Set-TextValue $ws.Range("E28") "  -3.96%  "

Set-TextValue $ws.Range("E29") "  -0.11%  "

Set-TextValue $ws.Range("D30") "3.910.81"
Set-TextValue $ws.Range("E30") "  -0.56%  "

Set-TextValue $ws.Range("E31") "  -4.83%  "

Set-TextValue $ws.Range("E32") "  -3.51%  "

Set-TextValue $ws.Range("D33") "30.03"
Set-TextValue $ws.Range("E33") "  -1.81%  "

Set-TextValue $ws.Range("E34") "  -3.18%  "

Set-TextValue $ws.Range("D35") "9.22"
Set-TextValue $ws.Range("E35") "  -0.53%  "

Set-TextValue $ws.Range("D37") "3.715.87"
Set-TextValue $ws.Range("E37") "  -0.70%  "

Set-TextValue $ws.Range("E38") "  -3.67%  "

Set-TextValue $ws.Range("E39") "  -8.33%  "

Set-TextValue $ws.Range("D40") "0.138"
Set-TextValue $ws.Range("E40") "  -1.05%  "

Set-TextValue $ws.Range("E41") "  -0.36%  "

Set-TextValue $ws.Range("D42") "5.80"
Set-TextValue $ws.Range("E42") "  -0.96%  "

Set-TextValue $ws.Range("E43") "  -0.03%  "

Set-TextValue $ws.Range("D45") "44.07"
Set-TextValue $ws.Range("E45") "  +9.06%  "

Set-TextValue $ws.Range("E46") "  -3.66%  "

Set-TextValue $ws.Range("D47") "46.81"
Set-TextValue $ws.Range("E47") "  +2.53%  "

Set-TextValue $ws.Range("E48") "  -1.80%  "

# Rows 50 and 51 swapped places (Monero <-> Bittensor) with refreshed prices
Set-TextValue $ws.Range("B50") "Bittensor"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D50") "390.31"
Set-TextValue $ws.Range("E50") "  -3.05%  "

Set-TextValue $ws.Range("B51") "Monero"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D51") "145.24"
Set-TextValue $ws.Range("E51") "  +2.24%  "
